$d = $word.ActiveDocument

function Clear-HighlightByFind($searchText) {
    $r = $d.Content
    $found = $r.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $r.HighlightColorIndex = 0
    }
    return $found
}

# 1) "או קנייה מעגלת הקניות (מנויים בלבד)." -> "או קנייה מעגלת הקניות." and drop the highlight
$r = $d.Content
$r.Find.Execute("או קנייה מעגלת הקניות (מנויים בלבד).", $false, $false, $false, $false, $false, $true, 1, $false, "או קנייה מעגלת הקניות.", 2)
Clear-HighlightByFind("או קנייה מעגלת הקניות.")

# 2) Merge the three runs about discount/policy rules support into one run with combined text
$r = $d.Content
$r.Find.Execute("הוספה תמיכה בהחלת חוקי הנחות ומדיניות הקניה במהלך תהליך רכישה של משתמש.", $false, $false, $false, $false, $false, $true, 1, $false, "הוספה תמיכה בהחלת חוקי הנחות ומדיניות הקניה במהלך תהליך רכישה של משתמש.", 2)

# 3) Move the _GoBack bookmark from after "watchtrans" to right after
#    "הוספת תמיכה בקניה מכמה חנויות ולא רק אחת...."
$r = $d.Content
$found = $r.Find.Execute("הוספת תמיכה בקניה מכמה חנויות ולא רק אחת....", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $r)
}

# 4) Drop highlight on "הוספנו עדכון של עגלת הקנייה לאחר קניה מוצלחת (היה חסר)."
Clear-HighlightByFind("הוספנו עדכון של עגלת הקנייה לאחר קניה מוצלחת (היה חסר).")

# 5) Drop highlight across the whole "לוודא שיורד ... / צפייה בעגלה" paragraph
#    (covers the inventory-rollback text, the mid-paragraph line break, and "צפייה בעגלה")
Clear-HighlightByFind("לוודא שיורד ממלאי החנות מה שנכשל חוזר למלאי החנות טרנזקציה מסתיימת מבוטלת")
Clear-HighlightByFind("צפייה בעגלה")

# 6) Drop highlight on "תיקון " and "תרחישי שימוש"
Clear-HighlightByFind("תיקון ")
Clear-HighlightByFind("תרחישי שימוש")

# 7) "תרחישי שימוש" is followed by ":" (highlighted) then " תהליך ה " -> re-partition into
#    ": תהליך" (no highlight) + " ה " (no highlight)
#    Do the second run first (strip the leading " תהליך" that is moving into the first run)
#    so the two runs never become formatting-identical before both edits have landed
#    (which would otherwise make the engine fold them back into a single run).
$rSetup = $d.Content
$foundSetup = $rSetup.Find.Execute("setup", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundSetup) {
    $rSetup.Collapse(1)
    $rSetup.MoveStart(1, -9)
    $rSetup.Text = " ה "
}

$r = $d.Content
$found = $r.Find.Execute("תרחישי שימוש", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r.Collapse(0)
    $r.MoveEnd(1, 1)
    $r.InsertAfter(" תהליך")
}

Clear-HighlightByFind(": תהליך")
